$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# Row 4: renewable energy consumption variable (name shortened)
$ws.Range("A4").Value = "renewable_eng"
$ws.Range("B4").Value = "Renewable energy consumption is the share of renewables energy in total final energy consumption."
$ws.Range("C4").Value = "numeric"
$ws.Range("D4").Value = "Renewable energy consumption (% of total final energy consumption)"

# Row 5: GDP per capita growth variable (name shortened + row height reduced)
$ws.Range("A5").Value = "gdppc_growth"
$ws.Range("B5").Value = "Annual percentage growth rate of GDP per capita based on constant local currency. Aggregates are based on constant 2010 U.S. dollars. GDP per capita is gross domestic product divided by midyear population. GDP at purchaser's prices is the sum of gross value added by all resident producers in the economy plus any product taxes and minus any subsidies not included in the value of the products. It is calculated without making deductions for depreciation of fabricated assets or for depletion and degradation of natural resources."
$ws.Range("C5").Value = "numeric"
$ws.Range("D5").Value = "GDP per capita growth (annual %)"
$ws.Rows.Item(5).RowHeight = 80

# Row 6: gross capital formation variable (name shortened)
$ws.Range("A6").Value = "capital_formation"
$ws.Range("B6").Value = 'Annual growth rate of gross capital formation based on constant local currency. Aggregates are based on constant 2010 U.S. dollars. Gross capital formation (formerly gross domestic investment) consists of outlays on additions to the fixed assets of the economy plus net changes in the level of inventories. Fixed assets include land improvements (fences, ditches, drains, and so on); plant, machinery, and equipment purchases; and the construction of roads, railways, and the like, including schools, offices, hospitals, private residential dwellings, and commercial and industrial buildings. Inventories are stocks of goods held by firms to meet temporary or unexpected fluctuations in production or sales, and "work in progress." According to the 1993 SNA, net acquisitions of valuables are also considered capital formation.'
$ws.Range("C6").Value = "numeric"
$ws.Range("D6").Value = "Gross capital formation (annual % growth)"

# Row 7: population growth variable (name shortened)
$ws.Range("A7").Value = "pop_growth"
$ws.Range("B7").Value = "Annual population growth rate for year t is the exponential rate of growth of midyear population from year t-1 to t, expressed as a percentage . Population is based on the de facto definition of population, which counts all residents regardless of legal status or citizenship."
$ws.Range("C7").Value = "numeric"
$ws.Range("D7").Value = "Population growth (annual %)"

# Row 8: replace "forest_area_percentage" variable with "trade"
$ws.Range("A8").Value = "trade"
$ws.Range("B8").Value = "Trade is the sum of exports and imports of goods and services measured as a share of gross domestic product."
$ws.Range("C8").Value = "numeric"
$ws.Range("D8").Value = "Trade (% of GDP)"
$ws.Rows.Item(8).AutoFit()

# Selection moves to B13 (next empty area below the table)
$ws.Range("B13").Select() | Out-Null
